# ALU function table maintenance:
#  - fix the wrong opcode that had been entered for operand "B"
#  - add two new rows for the "NOT A" / "NOT B" operations (inserted right
#    after the existing "A" / "B" rows, pushing NAND..CONST0 down by two)
#  - fix the opcode that had mistakenly been entered for NAND
#  - drop the old "RL" / "RR" rows, which are no longer part of the table
#  - refresh the selection to highlight the B column rows that were touched

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the old NAND row (row 11) to make room for
# "NOT A" and "NOT B". The new rows inherit the (highlighted) style of the
# row above them, same as Excel's native "Insert Copied/Shifted Cells".
$ws.Rows("11:12").Insert()

$ws.Range("A11").Value = "NOT A"
$ws.Range("A12").Value = "NOT B"
$ws.Range("B11").Value = "0b011100"

# Bug fix: operand "B" had the wrong opcode recorded.
$ws.Range("B10").Value = "0b011010"

# Bug fix: NAND's opcode (now on row 13 after the insert above) was wrong.
$ws.Range("B13").Value = "0b010011"

$ws.Range("B12").Value = "0b010101"

# The "RL" / "RR" rows (now at rows 21-22 after the insert above) are no
# longer needed; delete them and shift everything below back up.
$ws.Rows("21:22").Delete()

$ws.Range("B10:B13").Select() | Out-Null
